# Inserting data into a created table in a created SQLite database file
# Appends 5 new people (rows 5-9) read out of the (external) SQLite table
# into the "database" worksheet, in the same order the source data-entry
# form wrote them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Marina Osuna -> row 5 --------------------------------------------
$ws.Range("A5").Value = "Marina"
$ws.Range("B5").Value = "Osuna"
$ws.Range("C5").Value = "Mrs."
$ws.Range("D5").Value = "18"
$ws.Range("E5").Value = "United Kingdom"
$ws.Range("F5").Value = "woman"
$ws.Range("G5").Value = "registered"
$ws.Range("H5").Value = "3"
$ws.Range("I5").Value = "2"

# --- Fernando Domínguez -> row 6 (no title) ----------------------------
$ws.Range("A6").Value = "Fernando"
$ws.Range("B6").Value = "Domínguez"
$ws.Range("D6").Value = "25"
$ws.Range("E6").Value = "Spain"
$ws.Range("F6").Value = "man"
$ws.Range("G6").Value = "not registered"
$ws.Range("H6").Value = "0"
$ws.Range("I6").Value = "0"

# --- Lucrecia Hendrich -> row 9 -----------------------------------------
$ws.Range("A9").Value = "Lucrecia"
$ws.Range("B9").Value = "Hendrich"
$ws.Range("C9").Value = "MSc."
$ws.Range("D9").Value = "21"
$ws.Range("E9").Value = "Portugal"
$ws.Range("F9").Value = "woman"
$ws.Range("G9").Value = "registered"
$ws.Range("H9").Value = "13"
$ws.Range("I9").Value = "4"

# --- Ander Muñoz -> row 7 ------------------------------------------------
$ws.Range("A7").Value = "Ander"
$ws.Range("B7").Value = "Muñoz"
$ws.Range("C7").Value = "Dr."
$ws.Range("D7").Value = "19"
$ws.Range("E7").Value = "Greece"
$ws.Range("F7").Value = "man"
$ws.Range("G7").Value = "registered"
$ws.Range("H7").Value = "7"
$ws.Range("I7").Value = "3"

# --- Guzmán Osuna -> row 8 (no title) ------------------------------------
$ws.Range("A8").Value = "Guzmán"
$ws.Range("B8").Value = "Osuna"
$ws.Range("D8").Value = "20"
$ws.Range("E8").Value = "Spain"
$ws.Range("F8").Value = "man"
$ws.Range("G8").Value = "not registered"
$ws.Range("H8").Value = "3"
$ws.Range("I8").Value = "1"

# Column C ("Title") and column E ("Nationality") grow a bit wider to fit
# the newly-entered values ("Mrs."/"MSc." and "United Kingdom").
$ws.Range("C1").ColumnWidth = 4.1666666666
$ws.Range("E1").ColumnWidth = 13.1666666666

# Leave the cursor where the form's last write landed.
$ws.Range("M8").Select()
